# Commit: "add parquet, update storage with azcopy, more building"
#
# This edits the "class_schedule" sheet:
#  - C14 ("Thurs, Sep 24" / Big Data row): append a new bullet linking the
#    Parquet Format notebook to the "Do Before Class" list.
#  - C19 ("Tues, Oct 13" / Pandas Reshaping row): append a new bullet noting
#    the Project Strategy Plan is due, to the "Do Before Class" list.
#  - Row heights for rows 14 and 19 grow to fit the extra wrapped line.
#  - D19 picks up wrap-text formatting (cosmetic, matches the real edit).
#  - The active selection moves to C22, with the view scrolled toward A19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C14: append a Parquet Format bullet -----------------------------------
$c14 = $ws.Range("C14")
$c14Old = $c14.Value()
$c14.Value = $c14Old + "`n- ``Parquet Format <parquet.ipynb>```_"

# --- C19: append a Project Strategy Plan Due bullet -------------------------
$c19 = $ws.Range("C19")
$c19Old = $c19.Value()
$c19.Value = $c19Old + "`n- **Project Strategy Plan Due**"

# --- Row heights grow to fit the new wrapped line ---------------------------
$ws.Rows.Item(14).RowHeight = 68
$ws.Rows.Item(19).RowHeight = 99

# --- D19 picks up wrap text formatting --------------------------------------
$ws.Range("D19").WrapText = $true

# --- Update active selection / scroll position ------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("C22").Select()
